# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Rule row 11 (previously labelled "R40") is renamed to "1".
# Force the numeric-looking label to be stored as text (not a number),
# matching the original cell's string type.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
